# Update odds data in rows 5, 6, 7, 14, 15 and 17 of the FlashScore
# weekly-fixtures sheet to reflect the refreshed betting odds.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 2.5
$ws.Range("H5").Value = 2.8
$ws.Range("I5").Value = 3.4
$ws.Range("J5").Value = 3.5
$ws.Range("K5").Value = 1.8
$ws.Range("L5").Value = 4.33
$ws.Range("M5").Value = 1.14
$ws.Range("O5").Value = 1.67
$ws.Range("Q5").Value = 3.4
$ws.Range("R5").Value = 1.33
$ws.Range("T5").Value = 1.1
$ws.Range("W5").Value = 2.5
$ws.Range("X5").Value = 1.5
$ws.Range("Y5").Value = 5.5
$ws.Range("Z5").Value = 10
$ws.Range("AA5").Value = 11
$ws.Range("AB5").Value = 26
$ws.Range("AE5").Value = 5
$ws.Range("AI5").Value = 6.5
$ws.Range("AJ5").Value = 15
$ws.Range("AK5").Value = 15
$ws.Range("AL5").Value = 41
$ws.Range("AM5").Value = 41
$ws.Range("AR5").Value = 5.6
$ws.Range("AS5").Value = 1.15
$ws.Range("H6").Value = 3.2
$ws.Range("I6").Value = 4.5
$ws.Range("O6").Value = 1.5
$ws.Range("P6").Value = 2.5
$ws.Range("Q6").Value = 2.5
$ws.Range("R6").Value = 1.5
$ws.Range("S6").Value = 5
$ws.Range("T6").Value = 1.17
$ws.Range("U6").Value = 1.57
$ws.Range("V6").Value = 2.25
$ws.Range("W6").Value = 2.2
$ws.Range("X6").Value = 1.62
$ws.Range("AA6").Value = 9.5
$ws.Range("AE6").Value = 6.5
$ws.Range("AG6").Value = 21
$ws.Range("AH6").Value = 81
$ws.Range("AI6").Value = 9
$ws.Range("AP6").Value = 1.93
$ws.Range("AQ6").Value = 1.93
$ws.Range("AR6").Value = 4.1
$ws.Range("AS6").Value = 1.24
$ws.Range("I7").Value = 2.15
$ws.Range("J7").Value = 4.5
$ws.Range("L7").Value = 3.1
$ws.Range("Q7").Value = 3.4
$ws.Range("R7").Value = 1.33
$ws.Range("S7").Value = 7
$ws.Range("T7").Value = 1.1
$ws.Range("U7").Value = 1.7
$ws.Range("V7").Value = 2.08
$ws.Range("G14").Value = 2.75
$ws.Range("I14").Value = 3
$ws.Range("J14").Value = 3.6
$ws.Range("L14").Value = 3.75
$ws.Range("M14").Value = 1.13
$ws.Range("N14").Value = 6
$ws.Range("O14").Value = 1.57
$ws.Range("P14").Value = 2.25
$ws.Range("Q14").Value = 2.88
$ws.Range("R14").Value = 1.4
$ws.Range("W14").Value = 2.25
$ws.Range("X14").Value = 1.57
$ws.Range("Z14").Value = 11
$ws.Range("AJ14").Value = 13
$ws.Range("AK14").Value = 13
$ws.Range("AL14").Value = 34
$ws.Range("AR14").Value = 4.8
$ws.Range("G15").Value = 1.3
$ws.Range("H15").Value = 4.75
$ws.Range("I15").Value = 9
$ws.Range("J15").Value = 1.91
$ws.Range("K15").Value = 2.25
$ws.Range("L15").Value = 9.5
$ws.Range("M15").Value = 1.07
$ws.Range("N15").Value = 8.5
$ws.Range("O15").Value = 1.33
$ws.Range("P15").Value = 3.25
$ws.Range("Q15").Value = 2.08
$ws.Range("R15").Value = 1.73
$ws.Range("W15").Value = 2.63
$ws.Range("X15").Value = 1.44
$ws.Range("AB15").Value = 8
$ws.Range("AE15").Value = 8.5
$ws.Range("AF15").Value = 9.5
$ws.Range("AI15").Value = 17
$ws.Range("AK15").Value = 29
$ws.Range("AL15").Value = 126
$ws.Range("AM15").Value = 81
$ws.Range("AP15").Value = 1.58
$ws.Range("AQ15").Value = 2.34
$ws.Range("AR15").Value = 3.05
$ws.Range("AS15").Value = 1.37
$ws.Range("G17").Value = 2.45
$ws.Range("I17").Value = 2.55
$ws.Range("J17").Value = 3.6
$ws.Range("K17").Value = 2
$ws.Range("U17").Value = 1.5
$ws.Range("V17").Value = 2.5
$ws.Range("Y17").Value = 8
$ws.Range("AG17").Value = 15
$ws.Range("AN17").Value = 34
